# Add support for SEVA plasmids:
#  1. Insert a new "SEVASource" worksheet right after "WekWikGeneIdSource"
#     (a duplicate layout of WekWikGeneIdSource, with an extra "seva" option
#     in its repository_name dropdown list).
#  2. Append ",seva" to the repository_name/repository type dropdown list on
#     every existing sheet that already offered it.

$wb = $excel.ActiveWorkbook

$oldList = "addgene,genbank,benchling,snapgene,euroscarf,igem,wekwikgene"
$newList = "addgene,genbank,benchling,snapgene,euroscarf,igem,wekwikgene,seva"
$newListFormula = '"' + $newList + '"'

# --- 1. Create the new SEVASource sheet, positioned right after WekWikGeneIdSource ---

$afterSheet = $wb.Worksheets.Item("WekWikGeneIdSource")
$sevaSheet = $wb.Worksheets.Add($null, $afterSheet)
$sevaSheet.Name = "SEVASource"

$sevaSheet.Range("A1").Value = "sequence_file_url"
$sevaSheet.Range("B1").Value = "repository_id"
$sevaSheet.Range("C1").Value = "repository_name"
$sevaSheet.Range("D1").Value = "input"
$sevaSheet.Range("E1").Value = "output"
$sevaSheet.Range("F1").Value = "type"
$sevaSheet.Range("G1").Value = "output_name"
$sevaSheet.Range("H1").Value = "id"

$sevaSheet.Range("C2:C1048576").Validation.Add(3, 1, 1, $newListFormula)

# --- 2. Extend the existing repository-name dropdown lists with ",seva" ---

$sheetsAndRanges = @(
    @{ Sheet = "RepositoryIdSource";      Range = "B2:B1048576" },
    @{ Sheet = "AddGeneIdSource";         Range = "D2:D1048576" },
    @{ Sheet = "WekWikGeneIdSource";      Range = "C2:C1048576" },
    @{ Sheet = "BenchlingUrlSource";      Range = "B2:B1048576" },
    @{ Sheet = "SnapGenePlasmidSource";   Range = "B2:B1048576" },
    @{ Sheet = "EuroscarfSource";         Range = "B2:B1048576" },
    @{ Sheet = "IGEMSource";              Range = "C2:C1048576" }
)

foreach ($entry in $sheetsAndRanges) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $rng = $ws.Range($entry.Range)
    $rng.Validation.Formula1 = $newListFormula
}

Write-Host "SEVA source sheet added and dropdown lists updated"
